$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 255-267 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A255").NumberFormat = "@"
$ws.Range("A255").Value = "2026-01-28"
$ws.Range("B255").Value = "15:13:21"
$ws.Range("C255").Value = "15:00"
$ws.Range("D255").Value = "Bathroom"
$ws.Range("E255").Value = "No Motion"
$ws.Range("F255").Value = "Inactive"
$ws.Range("A256").NumberFormat = "@"
$ws.Range("A256").Value = "2026-01-28"
$ws.Range("B256").Value = "15:13:24"
$ws.Range("C256").Value = "15:00"
$ws.Range("D256").Value = "Bathroom"
$ws.Range("E256").Value = "No Motion"
$ws.Range("F256").Value = "Inactive"
$ws.Range("A257").NumberFormat = "@"
$ws.Range("A257").Value = "2026-01-28"
$ws.Range("B257").Value = "15:13:28"
$ws.Range("C257").Value = "15:00"
$ws.Range("D257").Value = "Bathroom"
$ws.Range("E257").Value = "No Motion"
$ws.Range("F257").Value = "Inactive"
$ws.Range("A258").NumberFormat = "@"
$ws.Range("A258").Value = "2026-01-28"
$ws.Range("B258").Value = "15:13:32"
$ws.Range("C258").Value = "15:00"
$ws.Range("D258").Value = "Bathroom"
$ws.Range("E258").Value = "No Motion"
$ws.Range("F258").Value = "Inactive"
$ws.Range("A259").NumberFormat = "@"
$ws.Range("A259").Value = "2026-01-28"
$ws.Range("B259").Value = "15:13:37"
$ws.Range("C259").Value = "15:00"
$ws.Range("D259").Value = "Bathroom"
$ws.Range("E259").Value = "No Motion"
$ws.Range("F259").Value = "Inactive"
$ws.Range("A260").NumberFormat = "@"
$ws.Range("A260").Value = "2026-01-28"
$ws.Range("B260").Value = "15:13:42"
$ws.Range("C260").Value = "15:00"
$ws.Range("D260").Value = "Bathroom"
$ws.Range("E260").Value = "No Motion"
$ws.Range("F260").Value = "Inactive"
$ws.Range("A261").NumberFormat = "@"
$ws.Range("A261").Value = "2026-01-28"
$ws.Range("B261").Value = "15:13:48"
$ws.Range("C261").Value = "15:00"
$ws.Range("D261").Value = "Bathroom"
$ws.Range("E261").Value = "No Motion"
$ws.Range("F261").Value = "Inactive"
$ws.Range("A262").NumberFormat = "@"
$ws.Range("A262").Value = "2026-01-28"
$ws.Range("B262").Value = "15:13:53"
$ws.Range("C262").Value = "15:00"
$ws.Range("D262").Value = "Bathroom"
$ws.Range("E262").Value = "No Motion"
$ws.Range("F262").Value = "Inactive"
$ws.Range("A263").NumberFormat = "@"
$ws.Range("A263").Value = "2026-01-28"
$ws.Range("B263").Value = "15:13:57"
$ws.Range("C263").Value = "15:00"
$ws.Range("D263").Value = "Bathroom"
$ws.Range("E263").Value = "No Motion"
$ws.Range("F263").Value = "Inactive"
$ws.Range("A264").NumberFormat = "@"
$ws.Range("A264").Value = "2026-01-28"
$ws.Range("B264").Value = "15:14:02"
$ws.Range("C264").Value = "15:00"
$ws.Range("D264").Value = "Bathroom"
$ws.Range("E264").Value = "No Motion"
$ws.Range("F264").Value = "Inactive"
$ws.Range("A265").NumberFormat = "@"
$ws.Range("A265").Value = "2026-01-28"
$ws.Range("B265").Value = "15:14:09"
$ws.Range("C265").Value = "15:00"
$ws.Range("D265").Value = "Bathroom"
$ws.Range("E265").Value = "No Motion"
$ws.Range("F265").Value = "Inactive"
$ws.Range("A266").NumberFormat = "@"
$ws.Range("A266").Value = "2026-01-28"
$ws.Range("B266").Value = "15:14:13"
$ws.Range("C266").Value = "15:00"
$ws.Range("D266").Value = "Bathroom"
$ws.Range("E266").Value = "No Motion"
$ws.Range("F266").Value = "Inactive"
$ws.Range("A267").NumberFormat = "@"
$ws.Range("A267").Value = "2026-01-28"
$ws.Range("B267").Value = "15:14:18"
$ws.Range("C267").Value = "15:00"
$ws.Range("D267").Value = "Bathroom"
$ws.Range("E267").Value = "No Motion"
$ws.Range("F267").Value = "Inactive"

# --- Humidity sheet: append rows 245-256 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A245").NumberFormat = "@"
$ws.Range("E245").NumberFormat = "@"
$ws.Range("A245").Value = "2026-01-28"
$ws.Range("B245").Value = "15:13:22"
$ws.Range("C245").Value = "15:00"
$ws.Range("D245").Value = "Bathroom"
$ws.Range("E245").Value = "88.4%"
$ws.Range("F245").Value = "Active"
$ws.Range("A246").NumberFormat = "@"
$ws.Range("E246").NumberFormat = "@"
$ws.Range("A246").Value = "2026-01-28"
$ws.Range("B246").Value = "15:13:27"
$ws.Range("C246").Value = "15:00"
$ws.Range("D246").Value = "Bathroom"
$ws.Range("E246").Value = "88.4%"
$ws.Range("F246").Value = "Active"
$ws.Range("A247").NumberFormat = "@"
$ws.Range("E247").NumberFormat = "@"
$ws.Range("A247").Value = "2026-01-28"
$ws.Range("B247").Value = "15:13:31"
$ws.Range("C247").Value = "15:00"
$ws.Range("D247").Value = "Bathroom"
$ws.Range("E247").Value = "88.5%"
$ws.Range("F247").Value = "Active"
$ws.Range("A248").NumberFormat = "@"
$ws.Range("E248").NumberFormat = "@"
$ws.Range("A248").Value = "2026-01-28"
$ws.Range("B248").Value = "15:13:35"
$ws.Range("C248").Value = "15:00"
$ws.Range("D248").Value = "Bathroom"
$ws.Range("E248").Value = "87.5%"
$ws.Range("F248").Value = "Active"
$ws.Range("A249").NumberFormat = "@"
$ws.Range("E249").NumberFormat = "@"
$ws.Range("A249").Value = "2026-01-28"
$ws.Range("B249").Value = "15:13:39"
$ws.Range("C249").Value = "15:00"
$ws.Range("D249").Value = "Bathroom"
$ws.Range("E249").Value = "88.4%"
$ws.Range("F249").Value = "Active"
$ws.Range("A250").NumberFormat = "@"
$ws.Range("E250").NumberFormat = "@"
$ws.Range("A250").Value = "2026-01-28"
$ws.Range("B250").Value = "15:13:47"
$ws.Range("C250").Value = "15:00"
$ws.Range("D250").Value = "Bathroom"
$ws.Range("E250").Value = "87.5%"
$ws.Range("F250").Value = "Active"
$ws.Range("A251").NumberFormat = "@"
$ws.Range("E251").NumberFormat = "@"
$ws.Range("A251").Value = "2026-01-28"
$ws.Range("B251").Value = "15:13:55"
$ws.Range("C251").Value = "15:00"
$ws.Range("D251").Value = "Bathroom"
$ws.Range("E251").Value = "87.5%"
$ws.Range("F251").Value = "Active"
$ws.Range("A252").NumberFormat = "@"
$ws.Range("E252").NumberFormat = "@"
$ws.Range("A252").Value = "2026-01-28"
$ws.Range("B252").Value = "15:13:59"
$ws.Range("C252").Value = "15:00"
$ws.Range("D252").Value = "Bathroom"
$ws.Range("E252").Value = "88.4%"
$ws.Range("F252").Value = "Active"
$ws.Range("A253").NumberFormat = "@"
$ws.Range("E253").NumberFormat = "@"
$ws.Range("A253").Value = "2026-01-28"
$ws.Range("B253").Value = "15:14:07"
$ws.Range("C253").Value = "15:00"
$ws.Range("D253").Value = "Bathroom"
$ws.Range("E253").Value = "88.4%"
$ws.Range("F253").Value = "Active"
$ws.Range("A254").NumberFormat = "@"
$ws.Range("E254").NumberFormat = "@"
$ws.Range("A254").Value = "2026-01-28"
$ws.Range("B254").Value = "15:14:11"
$ws.Range("C254").Value = "15:00"
$ws.Range("D254").Value = "Bathroom"
$ws.Range("E254").Value = "88.4%"
$ws.Range("F254").Value = "Active"
$ws.Range("A255").NumberFormat = "@"
$ws.Range("E255").NumberFormat = "@"
$ws.Range("A255").Value = "2026-01-28"
$ws.Range("B255").Value = "15:14:15"
$ws.Range("C255").Value = "15:00"
$ws.Range("D255").Value = "Bathroom"
$ws.Range("E255").Value = "87.4%"
$ws.Range("F255").Value = "Active"
$ws.Range("A256").NumberFormat = "@"
$ws.Range("E256").NumberFormat = "@"
$ws.Range("A256").Value = "2026-01-28"
$ws.Range("B256").Value = "15:14:19"
$ws.Range("C256").Value = "15:00"
$ws.Range("D256").Value = "Bathroom"
$ws.Range("E256").Value = "88.4%"
$ws.Range("F256").Value = "Active"

# --- Temperature sheet: append rows 245-256 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A245").NumberFormat = "@"
$ws.Range("A245").Value = "2026-01-28"
$ws.Range("B245").Value = "15:13:23"
$ws.Range("C245").Value = "15:00"
$ws.Range("D245").Value = "Bathroom"
$ws.Range("E245").Value = "23.0C"
$ws.Range("F245").Value = "Active"
$ws.Range("A246").NumberFormat = "@"
$ws.Range("A246").Value = "2026-01-28"
$ws.Range("B246").Value = "15:13:28"
$ws.Range("C246").Value = "15:00"
$ws.Range("D246").Value = "Bathroom"
$ws.Range("E246").Value = "23.0C"
$ws.Range("F246").Value = "Active"
$ws.Range("A247").NumberFormat = "@"
$ws.Range("A247").Value = "2026-01-28"
$ws.Range("B247").Value = "15:13:32"
$ws.Range("C247").Value = "15:00"
$ws.Range("D247").Value = "Bathroom"
$ws.Range("E247").Value = "23.0C"
$ws.Range("F247").Value = "Active"
$ws.Range("A248").NumberFormat = "@"
$ws.Range("A248").Value = "2026-01-28"
$ws.Range("B248").Value = "15:13:36"
$ws.Range("C248").Value = "15:00"
$ws.Range("D248").Value = "Bathroom"
$ws.Range("E248").Value = "22.9C"
$ws.Range("F248").Value = "Active"
$ws.Range("A249").NumberFormat = "@"
$ws.Range("A249").Value = "2026-01-28"
$ws.Range("B249").Value = "15:13:40"
$ws.Range("C249").Value = "15:00"
$ws.Range("D249").Value = "Bathroom"
$ws.Range("E249").Value = "22.9C"
$ws.Range("F249").Value = "Active"
$ws.Range("A250").NumberFormat = "@"
$ws.Range("A250").Value = "2026-01-28"
$ws.Range("B250").Value = "15:13:48"
$ws.Range("C250").Value = "15:00"
$ws.Range("D250").Value = "Bathroom"
$ws.Range("E250").Value = "23.0C"
$ws.Range("F250").Value = "Active"
$ws.Range("A251").NumberFormat = "@"
$ws.Range("A251").Value = "2026-01-28"
$ws.Range("B251").Value = "15:13:56"
$ws.Range("C251").Value = "15:00"
$ws.Range("D251").Value = "Bathroom"
$ws.Range("E251").Value = "22.9C"
$ws.Range("F251").Value = "Active"
$ws.Range("A252").NumberFormat = "@"
$ws.Range("A252").Value = "2026-01-28"
$ws.Range("B252").Value = "15:14:00"
$ws.Range("C252").Value = "15:00"
$ws.Range("D252").Value = "Bathroom"
$ws.Range("E252").Value = "22.9C"
$ws.Range("F252").Value = "Active"
$ws.Range("A253").NumberFormat = "@"
$ws.Range("A253").Value = "2026-01-28"
$ws.Range("B253").Value = "15:14:08"
$ws.Range("C253").Value = "15:00"
$ws.Range("D253").Value = "Bathroom"
$ws.Range("E253").Value = "22.9C"
$ws.Range("F253").Value = "Active"
$ws.Range("A254").NumberFormat = "@"
$ws.Range("A254").Value = "2026-01-28"
$ws.Range("B254").Value = "15:14:12"
$ws.Range("C254").Value = "15:00"
$ws.Range("D254").Value = "Bathroom"
$ws.Range("E254").Value = "22.9C"
$ws.Range("F254").Value = "Active"
$ws.Range("A255").NumberFormat = "@"
$ws.Range("A255").Value = "2026-01-28"
$ws.Range("B255").Value = "15:14:16"
$ws.Range("C255").Value = "15:00"
$ws.Range("D255").Value = "Bathroom"
$ws.Range("E255").Value = "22.9C"
$ws.Range("F255").Value = "Active"
$ws.Range("A256").NumberFormat = "@"
$ws.Range("A256").Value = "2026-01-28"
$ws.Range("B256").Value = "15:14:20"
$ws.Range("C256").Value = "15:00"
$ws.Range("D256").Value = "Bathroom"
$ws.Range("E256").Value = "23.0C"
$ws.Range("F256").Value = "Active"

# --- Proximity sheet: append rows 33-35 ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "2026-01-28"
$ws.Range("B33").Value = "15:13:22"
$ws.Range("C33").Value = "15:00"
$ws.Range("D33").Value = "Bedroom Door"
$ws.Range("E33").Value = "Detected"
$ws.Range("F33").Value = "Active"
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "2026-01-28"
$ws.Range("B34").Value = "15:13:24"
$ws.Range("C34").Value = "15:00"
$ws.Range("D34").Value = "Bedroom Door"
$ws.Range("E34").Value = "Clear"
$ws.Range("F34").Value = "Inactive"
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "2026-01-28"
$ws.Range("B35").Value = "15:14:04"
$ws.Range("C35").Value = "15:00"
$ws.Range("D35").Value = "Living Room Main Door"
$ws.Range("E35").Value = "ENTER"
$ws.Range("F35").Value = "User ENTERED Living Room Main Door"

# --- Camera sheet: append row 15 ---
$ws = $wb.Worksheets.Item("Camera")
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2026-01-28"
$ws.Range("B15").Value = "15:14:05"
$ws.Range("C15").Value = "15:00"
$ws.Range("D15").Value = "Living Room Main Door"
$ws.Range("E15").Value = "Image Captured"
$ws.Range("F15").Value = "Active"
